$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row above the current header row (row 1). Everything that
#    was row 1..7 shifts down to row 2..8.
# ---------------------------------------------------------------------------
$ws.Range("A1:L1").Insert(-4121)  # xlShiftDown

# ---------------------------------------------------------------------------
# 2) Build the new row 1 (the "OLS" group-header row above the real header).
#    Style it to match the bordered, non-bold look used throughout the table
#    by copying the format from the (now shifted) former J1 cell (J2), then
#    centering the A/D/E cells horizontally.
# ---------------------------------------------------------------------------
$ws.Range("J2").Copy()
foreach ($addr in @("A1","B1","C1","D1","E1","F1","G1","H1","I1","J1")) {
  $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
foreach ($addr in @("A1","D1","E1")) {
  $ws.Range($addr).HorizontalAlignment = -4108  # xlCenter
}
$excel.CutCopyMode = 0

# Merge D1:E1 and set its text to the new shared string "OLS".
$ws.Range("D1:E1").Merge()
$ws.Range("D1").Value = "OLS"

# ---------------------------------------------------------------------------
# 3) Update the real header row (now row 2). Un-bold A2:E2 (they keep their
#    border/alignment, just switch from the bold font to the regular one),
#    and set the two new header labels in D2/E2.
# ---------------------------------------------------------------------------
foreach ($addr in @("A2","B2","C2","D2","E2")) {
  $ws.Range($addr).Font.Bold = $false
}
$ws.Range("D2").Value = "(1) Market Value (euros)"
$ws.Range("E2").Value = "(2) Natural Log of Market Value"

# ---------------------------------------------------------------------------
# 4) Fill in the new data columns D/E for the six metric rows (now rows 3-8)
#    by copying the number format from column C (which already matches the
#    desired "font1 border/no-border + centered" style for each row) and
#    then writing the new values.
# ---------------------------------------------------------------------------
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D3").Value = 6426873.3690999998
$ws.Range("E3").Value = 0.38640000000000002

$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D4").Value = 106516475775883
$ws.Range("E4").Value = 0.25890000000000002

$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D5").Value = 10320681.943399999
$ws.Range("E5").Value = 0.50880000000000003

$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D6").Value = 1.7145999999999999
$ws.Range("E6").Value = 0.024899999999999999

$ws.Range("C7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D7").Value = 4481399.3886000002
$ws.Range("E7").Value = 0.31090000000000001

# Row 8 (R-Squared) uses the percent-styled C8 as the template for D8/E8.
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D8").Value = 0.81899999999999995
$ws.Range("E8").Value = 0.878

# ---------------------------------------------------------------------------
# 5) Column widths: D grows (loses "best fit"), and the new column E is
#    added with its own best-fit width.
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 34.42578125
$ws.Columns("E").ColumnWidth = 36

# ---------------------------------------------------------------------------
# 6) Cosmetic bits that mirror the rest of the diff.
# ---------------------------------------------------------------------------
$ws.Range("C14").Select()
